# Oregon Samples TP 21 12/11/2019
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 73, 74, 76 previously carried the "Crm opened 11/19/2020" note (shared
# string that is being re-purposed below); they now reference the existing
# "Crm opened 11/19/2019" text instead.
$ws.Range("F73").Value = "Crm opened 11/19/2019"
$ws.Range("F74").Value = "Crm opened 11/19/2019"
$ws.Range("F76").Value = "Crm opened 11/19/2019"

# New data row 77 for the Oregon sample taken 12/11/2019.
$ws.Range("A77").Value = (Get-Date -Year 2019 -Month 12 -Day 11)
$ws.Range("B77").Value = 2208.24007259168
$ws.Range("C77").Value = 2207.0300000000002
$ws.Range("D77").Formula = "=100*(B77-C77)/C77"
$ws.Range("E77").Value = 169
$ws.Range("F77").Value = "New CRM opened 12/11/2019"

$ws.Range("D78").Select()
